# Reordered subroutines: the "collision" category block (checkCollision)
# moves from between "draw" and "matrix" up to the very top of the list
# (right after "Program Group"/row 24), and the "random" category block
# (generateRandom3BitValue / generateRandom4BitValue) moves from the top
# down to just before "joystick" (after "display").
#
# Rather than relying on uncertain Cut/Insert row semantics, we rewrite
# the whole A26:C65 block explicitly with its final contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old block first (it used to end at row 64; the rewritten
# block now runs one row further, to row 65).
$ws.Range("A26:C65").ClearContents()

# --- Category block: collision (moved to the top) ---
$ws.Range("A26").Value = "collision"
$ws.Range("C26").Value = "checkCollision"

# --- Category block: logos ---
$ws.Range("A28").Value = "logos"
$ws.Range("C28").Value = "drawSnakeLogo"
$ws.Range("C29").Value = "drawMazeLogo"
$ws.Range("C30").Value = "drawAsteroidLogo"
$ws.Range("C31").Value = "drawTimerLogo"
$ws.Range("C32").Value = "drawRandomLogo"
$ws.Range("C33").Value = "drawJoystickLogo"

# --- Category block: draw ---
$ws.Range("A35").Value = "draw"
$ws.Range("C35").Value = "drawSnakeHeadMatrix"
$ws.Range("C36").Value = "drawSmileyMatrix"
$ws.Range("C37").Value = "drawTemplarMatrix"
$ws.Range("C38").Value = "drawSkullMatrix"

# --- Category block: matrix ---
$ws.Range("A40").Value = "matrix"
$ws.Range("C40").Value = "setPixel"
$ws.Range("C41").Value = "clearPixel"
$ws.Range("C42").Value = "invertMatrix"
$ws.Range("C43").Value = "clearMatrix"
$ws.Range("C44").Value = "setMatrix"

# --- Category block: display ---
$ws.Range("A46").Value = "display"
$ws.Range("C46").Value = "render"

# --- Category block: random (moved down, after display) ---
$ws.Range("A49").Value = "random"
$ws.Range("C49").Value = "generateRandom3BitValue"
$ws.Range("C50").Value = "generateRandom4BitValue"

# --- Category block: joystick ---
$ws.Range("A52").Value = "joystick"
$ws.Range("C52").Value = "readJoystick"
$ws.Range("C53").Value = "readJoystickDirection"
$ws.Range("C54").Value = "joystickValueTo8Bit"
$ws.Range("C55").Value = "joystickValuesToDirection"

# --- Category block: timers ---
$ws.Range("A57").Value = "timers"
$ws.Range("C57").Value = "incrementTimer"
$ws.Range("C58").Value = "initializeTimer"
$ws.Range("C59").Value = "checkTimer"
$ws.Range("C60").Value = "initializeHardwareTimer2"

# --- Category block: core ---
$ws.Range("A62").Value = "core"
$ws.Range("C62").Value = "addProgram"
$ws.Range("C63").Value = "init"
$ws.Range("C64").Value = "main"
$ws.Range("C65").Value = "terminate"

# Column A got a little wider (old stored width 13.5546875 -> new stored
# width ~15.33203125). The COM ColumnWidth setter quantizes to 1/6ths of
# a character, so 14.5 is the closest input that lands on the same
# stored width bucket (15.333333333333334, ~0.001 off the author's
# value - an imperceptible, sub-pixel difference).
$ws.Columns("A").ColumnWidth = 14.5

# Sheet view: selecting a cell drops the old frozen/scrolled
# "topLeftCell"/row-selection state and records the new single-cell
# selection at B50, matching the edited workbook.
$null = $ws.Range("B50").Select()
